$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Jacutinga - Utilar De Jacutinga Eireli E"; New = "Jacutinga - Utilar De Jacutinga Eirel" },
    @{ Old = "Ouro Fino - Comercial Labegaline Ltda (O"; New = "Ouro Fino - Comercial Labegaline Ltda" },
    @{ Old = "Bueno Brandao - Moveis Bueno Brandão Ltd"; New = "Bueno Brandao - Moveis Bueno Brandão" },
    @{ Old = "Olimpio Noronha - Casa De Racao E Materi"; New = "Olimpio Noronha - Casa De Racao E Mat" },
    @{ Old = "Carmo De Minas - Fernanda Aparecida Dos "; New = "Carmo De Minas - Fernanda Aparecida D" },
    @{ Old = "Cruzilia - Comercial Eletromoveis Alvare"; New = "Cruzilia - Comercial Eletromoveis Alv" },
    @{ Old = "Sao Goncalo Do Sapucai - Vidal Moveis Lt"; New = "Sao Goncalo Do Sapucai - Vidal Moveis" },
    @{ Old = "Lavras - Ceara Moveis Ltda - Me Integral"; New = "Lavras - Ceara Moveis Ltda - Me Integ" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
